$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Preferences Satisfaction")

$ws.Range("E2").Value = "1, 2, 4, 7, 8, 9"
$ws.Range("F2").Value = "1, 2, 3, 4, 5, 7"
$ws.Range("H2").Value = 83.33333333333334
$ws.Range("J2").Value = 97.01
$ws.Range("J3").Value = 2.99
$ws.Range("E6").Value = "4, 6, 7, 8, 10"
$ws.Range("F6").Value = "1, 5, 2, 7, 3"
$ws.Range("H6").Value = 80
$ws.Range("E7").Value = "1, 4, 7, 8, 9"
$ws.Range("F7").Value = "4, 3, 2, 1, 7"
$ws.Range("H7").Value = 80
$ws.Range("E9").Value = "'8"
$ws.Range("F9").Value = "'2"
$ws.Range("E10").Value = "'9"
$ws.Range("F10").Value = "'2"
$ws.Range("E11").Value = "2, 4, 5, 6, 9, 10"
$ws.Range("F11").Value = "2, 1, 3, 4, 7, 5"
$ws.Range("H11").Value = 83.33333333333334
$ws.Range("E12").Value = "7, 8, 9, 10"
$ws.Range("F12").Value = "1, 6, 2, 4"
$ws.Range("E14").Value = "'9"
$ws.Range("F14").Value = "'2"
$ws.Range("E16").Value = "4, 5, 6, 8, 9, 10"
$ws.Range("F16").Value = "3, 5, 4, 8, 2, 1"
$ws.Range("H16").Value = 83.33333333333334
$ws.Range("E17").Value = "1, 4, 8"
$ws.Range("F17").Value = "3, 2, 5"
$ws.Range("E18").Value = "3, 5"
$ws.Range("F18").Value = "1, 3"
$ws.Range("E20").Value = "1, 2, 4, 6"
$ws.Range("F20").Value = "2, 1, 3, 5"
$ws.Range("E22").Value = "3, 5, 9, 10"
$ws.Range("F22").Value = "1, 6, 5, 2"
$ws.Range("E25").Value = "3, 5"
$ws.Range("F25").Value = "1, 3"
$ws.Range("E27").Value = "1, 6, 7"
$ws.Range("F27").Value = "1, 5, 3"
$ws.Range("E34").Value = "3, 9"
$ws.Range("F34").Value = "2, 4"
$ws.Range("E39").Value = "1, 2, 4, 6, 9, 10"
$ws.Range("F39").Value = "3, 1, 2, 4, 8, 5"
$ws.Range("H39").Value = 83.33333333333334
$ws.Range("E42").Value = "2, 3, 4, 6, 8, 10"
$ws.Range("F42").Value = "2, 3, 1, 4, 7, 5"
$ws.Range("H42").Value = 83.33333333333334
$ws.Range("E43").Value = "4, 5, 6, 10"
$ws.Range("F43").Value = "1, 2, 5, 4"
$ws.Range("E45").Value = "1, 2, 4, 7, 8, 9"
$ws.Range("F45").Value = "1, 2, 3, 5, 4, 7"
$ws.Range("H45").Value = 83.33333333333334
$ws.Range("E46").Value = "1, 3"
$ws.Range("F46").Value = "2, 3"
$ws.Range("E49").Value = "1, 2, 4, 5, 7, 8"
$ws.Range("F49").Value = "1, 2, 5, 7, 3, 4"
$ws.Range("H49").Value = 83.33333333333334
$ws.Range("E50").Value = "2, 6, 7, 9"
$ws.Range("F50").Value = "1, 3, 2, 6"
$ws.Range("E51").Value = "2, 4, 6, 9"
$ws.Range("F51").Value = "1, 2, 3, 6"
$ws.Range("E52").Value = "1, 3, 7, 8, 9"
$ws.Range("F52").Value = "2, 1, 5, 4, 7"
$ws.Range("H52").Value = 80
$ws.Range("E54").Value = "6, 9, 10"
$ws.Range("F54").Value = "4, 2, 1"
$ws.Range("E56").Value = "4, 6, 8, 10"
$ws.Range("F56").Value = "4, 3, 6, 1"
$ws.Range("E57").Value = "2, 5, 6, 9"
$ws.Range("F57").Value = "1, 3, 2, 5"
$ws.Range("E60").Value = "3, 8"
$ws.Range("F60").Value = "1, 4"
